# Generate Report for Handback
# Update handoff/handback timestamps for the file
# "1de47477-17a8-4f7c-aaae-b958c9749820.md" (row 3 in each sheet) to
# reflect a new handback report generation cycle.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 3: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$wsZhCn.Range("H3").Value = "2016-10-27 08:15:51"
$wsZhCn.Range("K3").Value = "2016-10-27 08:16:41"

# de-de sheet, row 3: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$wsDeDe.Range("H3").Value = "2016-10-27 08:16:04"
$wsDeDe.Range("K3").Value = "2016-10-27 08:16:59"

# Overview sheet, row 3: Latest HO Xliff Generate Date (G)
$wsOverview.Range("G3").Value = "2016-10-27 08:16:04"
